# Atualização de bases das ligas, do dia: 10-06-2024 às 21:53
#
# The source rows hold match-odds records. This edit re-sorts a handful of
# rows by re-distributing their B:AD payload (id/teams/odds/etc.) while the
# leading "A" column (sequential row number) stays put:
#   - rows 93-97 : cyclic shift (each row takes the payload that used to sit
#                  one row above it; row 93 wraps around and takes row 97's
#                  old payload)
#   - rows 100/101, 116/117, 173/174, 205/206 : simple pairwise swaps
#
# We read every source row's B:AD cells first (snapshotting with .Value2,
# which preserves the exact stored number/string — no literal retyping that
# could perturb floating point formatting), then write the snapshots back
# into their destination rows. Snapshotting first is required because the
# rotation is an in-place permutation with overlapping source/destination
# rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstCol = 2   # column B
$lastCol  = 30  # column AD

function Get-RowSnapshot($row) {
    $vals = @{}
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $vals[$c] = $ws.Cells.Item($row, $c).Value2
    }
    return $vals
}

function Set-RowSnapshot($row, $vals) {
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $ws.Cells.Item($row, $c).Value2 = $vals[$c]
    }
}

# --- snapshot every source row up front ---------------------------------
$snap93  = Get-RowSnapshot 93
$snap94  = Get-RowSnapshot 94
$snap95  = Get-RowSnapshot 95
$snap96  = Get-RowSnapshot 96
$snap97  = Get-RowSnapshot 97

$snap100 = Get-RowSnapshot 100
$snap101 = Get-RowSnapshot 101

$snap116 = Get-RowSnapshot 116
$snap117 = Get-RowSnapshot 117

$snap173 = Get-RowSnapshot 173
$snap174 = Get-RowSnapshot 174

$snap205 = Get-RowSnapshot 205
$snap206 = Get-RowSnapshot 206

# --- write the rotated / swapped payloads back ---------------------------
# rows 93-97: cyclic shift down by one (93 <- 97 <- 96 <- 95 <- 94 <- 93)
Set-RowSnapshot 93 $snap97
Set-RowSnapshot 94 $snap93
Set-RowSnapshot 95 $snap94
Set-RowSnapshot 96 $snap95
Set-RowSnapshot 97 $snap96

# simple pairwise swaps
Set-RowSnapshot 100 $snap101
Set-RowSnapshot 101 $snap100

Set-RowSnapshot 116 $snap117
Set-RowSnapshot 117 $snap116

Set-RowSnapshot 173 $snap174
Set-RowSnapshot 174 $snap173

Set-RowSnapshot 205 $snap206
Set-RowSnapshot 206 $snap205
